# MCU8APPS-52247: add a link back to the project's repository browse page
# on the last slide of the roadmap deck (below the existing "Next steps"
# graphic), as a plain auto-sized, non-wrapping text box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)

# Left/Top/Width/Height are expressed in points; the values below are the
# exact EMU target (50260,5855925)-(10278968x369332) converted at
# 914400 EMU/in = 12700 EMU/pt so the saved geometry lands on the same EMU.
$left   = 3.9574803150
$top    = 461.0964566929
$width  = 809.3675590551
$height = 29.0812598425

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 2"

$tb.TextFrame.TextRange.Text = "https://bitbucket.microchip.com/projects/MCU8NPIAPPS/repos/gsm_tracker_internship/browse"

# Auto-fit the box to its single line of text and don't wrap - matches the
# "type once, box grows to fit" text box PowerPoint creates by default.
$tb.TextFrame.WordWrap = 0
$tb.TextFrame.AutoSize = 1

# No background fill on the box.
$tb.Fill.Visible = 0
